$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (preserves text formatting, avoids numeric auto-conversion)
$updates = @{
    'D2' = '327.01'
    'E2' = '-1.09%'
    'D3' = '44.03'
    'E3' = '5.94%'
    'D4' = '5.478'
    'E4' = '-3.54%'
    'D5' = '0.08069'
    'E5' = '-4.26%'
    'D6' = '8.627'
    'E6' = '-2.06%'
    'D7' = '4.273'
    'E7' = '-4.72%'
    'D8' = '1.879'
    'E8' = '-5.52%'
    'D10' = '0.9372'
    'E10' = '1.04%'
    'D11' = '0.1167'
    'E11' = '-8.44%'
    'E12' = '-3.65%'
    'D13' = '0.09573'
    'E13' = '2.44%'
    'D14' = '0.04162'
    'E14' = '4.33%'
    'E15' = '0.32%'
    'D16' = '0.001267'
    'E16' = '-2.64%'
    'D17' = '0.006000'
    'E17' = '-1.84%'
    'D18' = '3.573'
    'E18' = '4.27%'
    'E19' = '-0.76%'
    'D20' = '8.580'
    'E20' = '-4.41%'
    'D21' = '0.1364'
    'E21' = '-0.09%'
    'D22' = '0.2594'
    'E22' = '3.21%'
    'D23' = '0.04352'
    'E23' = '-1.39%'
    'D24' = '0.001235'
    'E24' = '-0.70%'
    'D25' = '0.004374'
    'E25' = '-0.54%'
    'D26' = '0.0001234'
    'E26' = '3.57%'
    'D27' = '0.0004002'
    'E27' = '0.19%'
    'D39' = '0.02647'
    'E39' = '-7.55%'
    'D40' = '0.05442'
    'E40' = '-1.50%'
    'D41' = '0.01145'
    'E41' = '27.57%'
    'D42' = '0.007670'
    'E42' = '-2.84%'
    'D43' = '0.1391'
    'E43' = '-3.26%'
    'D44' = '0.002115'
    'E44' = '1.53%'
    'D45' = '0.009674'
    'E45' = '-11.78%'
    'D46' = '0.00006922'
    'E46' = '-4.49%'
    'D47' = '0.00000000753'
    'E47' = '0.19%'
    'D48' = '0.003569'
    'E48' = '10.69%'
    'D49' = '0.002278'
    'E49' = '-0.15%'
    'D50' = '0.00002107'
    'E50' = '0.19%'
    'D51' = '0.0002007'
    'E51' = '0.19%'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
